$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Map of cell address -> new text value. NumberFormat is forced to "@" (text)
# before assignment so that numeric-looking strings (e.g. "579.10", "1.00")
# are preserved verbatim instead of being coerced into numbers.
$changes = [ordered]@{
    "D2"  = "67.012.08"
    "D3"  = "3.235.92"
    "E3"  = "  +2.34%  "
    "E4"  = "  -0.02%  "
    "D5"  = "579.10"
    "E5"  = "  +4.09%  "
    "D6"  = "175.76"
    "E6"  = "  +3.28%  "
    "B7"  = "XRP"
    "C7"  = "https://coinranking.com/coin/-l8Mn2pVlRs-p+xrp-xrp"
    "D7"  = "0.606"
    "E7"  = "  +0.21%  "
    "B8"  = "USDC"
    "C8"  = "https://coinranking.com/coin/aKzUVe4Hh_CON+usdc-usdc"
    "D8"  = "1.00"
    "E8"  = "  -0.03%  "
    "D9"  = "3.234.31"
    "E9"  = "  +2.37%  "
    "E10" = "  +4.75%  "
    "D11" = "6.69"
    "E11" = "  +1.85%  "
    "E12" = "  +3.20%  "
    "D13" = "3.801.18"
    "E13" = "  +2.45%  "
    "E14" = "  +1.83%  "
    "D15" = "27.72"
    "E15" = "  +1.55%  "
    "D16" = "66.963.50"
    "E16" = "  +4.45%  "
    "E17" = "  +3.78%  "
    "D18" = "3.246.64"
    "E18" = "  +2.50%  "
    "E19" = "  +2.96%  "
    "D20" = "13.23"
    "E20" = "  +2.13%  "
    "D21" = "366.48"
    "E21" = "  +4.48%  "
    "D22" = "7.44"
    "E22" = "  +4.24%  "
    "E23" = "  +0.21%  "
    "D24" = "70.13"
    "E24" = "  +2.03%  "
    "D25" = "3.378.14"
    "E25" = "  +2.25%  "
    "E26" = "  +1.19%  "
    "E27" = "  +0.72%  "
    "E28" = "  +3.65%  "
    "E29" = "  +2.36%  "
    "E30" = "  +0.31%  "
    "D31" = "1.98"
    "E31" = "  +5.68%  "
    "D32" = "5.61"
    "E32" = "  +0.96%  "
    "E33" = "  +2.02%  "
    "E34" = "  -0.04%  "
    "D35" = "173.60"
    "E35" = "  +10.75%  "
    "E36" = "  +3.74%  "
    "D37" = "6.73"
    "E37" = "  +2.41%  "
    "E38" = "  +5.14%  "
    "E39" = "  +6.27%  "
    "E40" = "  +10.80%  "
    "D41" = "26.59"
    "E41" = "  +2.80%  "
    "D42" = "2.56"
    "E42" = "  +3.34%  "
    "D43" = "6.41"
    "E43" = "  +7.48%  "
    "D44" = "2.706.27"
    "E44" = "  +2.85%  "
    "D45" = "4.27"
    "E45" = "  +3.70%  "
    "E46" = "  +4.62%  "
    "D47" = "0.0669"
    "E47" = "  +3.41%  "
    "E48" = "  +5.12%  "
    "D49" = "333.79"
    "E49" = "  +3.07%  "
    "D50" = "0.0278"
    "E50" = "  +3.65%  "
    "D51" = "0.104"
    "E51" = "  +2.54%  "
}

foreach ($addr in $changes.Keys) {
    $rng = $ws.Range($addr)
    $rng.NumberFormat = "@"
    $rng.Value = $changes[$addr]
}
